$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.154.65'
$ws.Range("E2").Value = '  +3.16%  '

$ws.Range("D3").Value = '2.308.27'
$ws.Range("E3").Value = '  +2.10%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '310.62'
$ws.Range("E5").Value = '  +1.92%  '

$ws.Range("D6").Value = '101.23'
$ws.Range("E6").Value = '  +6.17%  '

$ws.Range("D7").Value = '0.538'
$ws.Range("E7").Value = '  +2.23%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  +7.28%  '

$ws.Range("D10").Value = '36.00'
$ws.Range("E10").Value = '  +2.87%  '

$ws.Range("E11").Value = '  +3.63%  '

$ws.Range("E12").Value = '  +1.06%  '

$ws.Range("E13").Value = '  +7.00%  '

$ws.Range("D14").Value = '2.665.93'
$ws.Range("E14").Value = '  +2.00%  '

$ws.Range("D15").Value = '15.03'
$ws.Range("E15").Value = '  +4.68%  '

$ws.Range("D16").Value = '2.301.32'
$ws.Range("E16").Value = '  +2.15%  '

$ws.Range("D17").Value = '0.815'
$ws.Range("E17").Value = '  +3.14%  '

$ws.Range("D18").Value = '43.104.54'
$ws.Range("E18").Value = '  +3.20%  '

$ws.Range("D19").Value = '12.55'
$ws.Range("E19").Value = '  +1.64%  '

$ws.Range("E20").Value = '  +2.42%  '

$ws.Range("E21").Value = '  +2.58%  '

$ws.Range("D22").Value = '68.57'
$ws.Range("E22").Value = '  +0.90%  '

$ws.Range("D23").Value = '241.06'
$ws.Range("E23").Value = '  +1.82%  '

$ws.Range("E24").Value = '  +4.60%  '

$ws.Range("E25").Value = '  +2.59%  '

$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.32%  '

$ws.Range("E27").Value = '  +5.13%  '

$ws.Range("D28").Value = '37.54'
$ws.Range("E28").Value = '  +2.51%  '

$ws.Range("D29").Value = '9.65'
$ws.Range("E29").Value = '  +2.32%  '

$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("D31").Value = '167.17'
$ws.Range("E31").Value = '  +4.41%  '

$ws.Range("E32").Value = '  +2.40%  '

$ws.Range("E33").Value = '  +0.11%  '

$ws.Range("E34").Value = '  -1.51%  '

$ws.Range("D35").Value = '17.72'
$ws.Range("E35").Value = '  +4.69%  '

$ws.Range("E36").Value = '  +1.10%  '

$ws.Range("E37").Value = '  +3.14%  '

$ws.Range("E38").Value = '  +0.91%  '

$ws.Range("E39").Value = '  +2.09%  '

$ws.Range("D40").Value = '1.84'
$ws.Range("E40").Value = '  +1.28%  '

$ws.Range("E41").Value = '  +8.09%  '

$ws.Range("E42").Value = '  -0.26%  '

$ws.Range("D43").Value = '1.988.29'
$ws.Range("E43").Value = '  +0.84%  '

$ws.Range("E44").Value = '  +2.79%  '

$ws.Range("D45").Value = '19.08'
$ws.Range("E45").Value = '  +1.58%  '

$ws.Range("E46").Value = '  +3.81%  '

$ws.Range("D47").Value = '9.85'
$ws.Range("E47").Value = '  -0.05%  '

$ws.Range("E48").Value = '  +18.35%  '

$ws.Range("D49").Value = '55.62'
$ws.Range("E49").Value = '  +5.42%  '

$ws.Range("D50").Value = '2.533.06'
$ws.Range("E50").Value = '  +1.85%  '

$ws.Range("E51").Value = '  +2.19%  '
